# Clear the previously-entered credentials/URL test data on the Input_Value
# sheet (L2: URL, M2: UserName, N2: Password) and remove the hyperlink that
# was attached to L2, then leave that range selected to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Remove the hyperlink that lived on L2 (pointed at the Oracle Cloud URL)
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# Clear the cell contents for L2:N2 (URL, UserName, Password test data)
$ws.Range("L2:N2").ClearContents()

# Reflect the final selection state (L2:N2 highlighted)
$ws.Range("L2:N2").Select()
